# Finished Field Download component.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlTop = -4160
$xlLeft = -4131

# --- 1. Rename the worksheet/tab. --------------------------------------
$ws.Name = "food_event_fields"

# --- 2. Rebuild each data row without the row-level "customFormat" flag
#        while preserving the exact per-cell formatting (bold header font,
#        top/left alignment, wrap text on the long description column) and
#        the explicit row heights already set on the wrapped rows. -------
$rowHeights = @{ 5 = 90; 7 = 135; 8 = 240; 9 = 240; 10 = 409; 11 = 270; 12 = 409 }

for ($r = 1; $r -le 12; $r++) {
    $ws.Rows.Item($r).ClearFormats()

    # column A & C: vertical-top only (bold only on the header row)
    $ws.Range("A$r").VerticalAlignment = $xlTop
    $ws.Range("C$r").VerticalAlignment = $xlTop

    # column B: left + vertical-top (bold only on the header row)
    $ws.Range("B$r").HorizontalAlignment = $xlLeft
    $ws.Range("B$r").VerticalAlignment = $xlTop

    # column D: vertical-top always; wrap text only on the rows that carry
    # an explicit (taller) row height
    $ws.Range("D$r").VerticalAlignment = $xlTop
    if ($rowHeights.ContainsKey($r)) {
        $ws.Range("D$r").WrapText = $true
    }

    if ($r -eq 1) {
        $ws.Range("A1:D1").Font.Bold = $true
    }

    if ($rowHeights.ContainsKey($r)) {
        $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
    }
}

# --- 3. Move the active selection from D5 to D8. ------------------------
$ws.Range("D8").Select()

# --- 4. Adjust the remembered workbook window position/size. -----------
$wb.Windows.Item(1).Left = 240
$wb.Windows.Item(1).Top = 240
$wb.Windows.Item(1).Width = 25360
$wb.Windows.Item(1).Height = 15820
